$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update C2:C6 values from 510 to 505
$ws.Range("C2:C6").Value = 505

# Update the active cell selection to C9
$ws.Range("C9").Select()
